# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.01 = 24461.88 pesos`n✅ 24461.88 pesos = 5.97 = 977.27 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 166.3
$wsTasas.Range("O10").Value = 4068.01
$wsTasas.Range("N12").Value = 4097.3
$wsTasas.Range("O12").Value = 163.69
